$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, pushing the existing rows 103:179 down to 104:180.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new price-report entry.
$ws.Cells.Item(103, 1).Value = 3
$ws.Cells.Item(103, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(103, 3).Value = "Coquimbo"
$ws.Cells.Item(103, 4).Value = 44673
$ws.Cells.Item(103, 5).Value = 5
$ws.Cells.Item(103, 6).Value = 100112030
$ws.Cells.Item(103, 7).Value = "Poroto granado"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 40
$ws.Cells.Item(103, 11).Value = 23000
$ws.Cells.Item(103, 12).Value = 23000
$ws.Cells.Item(103, 13).Value = 23000
$ws.Cells.Item(103, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(103, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(103, 16).Value = 920
$ws.Cells.Item(103, 17).Value = 25
$ws.Cells.Item(103, 18).Value = "Hortaliza"
